$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Mark most of the "checklist" rows with an "X" (same value already used in C1),
# matching the conceptual-model rendering that is now done in uml.
$ws.Range("C4").Value = "X"
$ws.Range("C5").Value = "X"
$ws.Range("C6").Value = "X"
$ws.Range("C9").Value = "X"
$ws.Range("C10").Value = "X"
$ws.Range("C11").Value = "X"

# Last row uses a distinct (new) shared string "x" (lowercase).
$ws.Range("C12").Value = "x"

# Update the remembered selection in the frozen (bottom-left) pane from K13 to C13.
$ws.Range("C13").Select()

$wb.Save()
